# Generate Report for Handback
# Updates the localization-status workbook: marks zh-cn/de-de handback rows
# as "Handed back: in sync with en-US", fills in Latest Target File / Latest
# Handback File / Latest Handback DateTime columns, links the new target-file
# cells to the same source doc on GitHub, and widens the columns that now
# hold longer text.

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78459990ababa072644d766c063513f7ef8462db/e2e/a7560759-331f-4243-83c1-b7e49ded50f9.md"
$urlE = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78459990ababa072644d766c063513f7ef8462db/e2e/e1f29d49-765f-4abe-8a5d-7a268dab63cd.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a7560759-331f-4243-83c1-b7e49ded50f9.md"
$wsZh.Range("J2").Value = "a7560759-331f-4243-83c1-b7e49ded50f9.b5f762ddd46277abdae7fc8847cc2bcde2ab6a0f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 07:47:57"

$wsZh.Range("I3").Value = "e1f29d49-765f-4abe-8a5d-7a268dab63cd.md"
$wsZh.Range("J3").Value = "e1f29d49-765f-4abe-8a5d-7a268dab63cd.f93c273936f3d89f087d4d48a10478ef87b4e57f.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 07:47:57"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a7560759-331f-4243-83c1-b7e49ded50f9.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlE, "", "", "e1f29d49-765f-4abe-8a5d-7a268dab63cd.md")

$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Range("I3").Style = "HyperLink"
$wsZh.Range("I2:I3").Font.Underline = $true
$wsZh.Range("I2:I3").Font.Color = 15570276

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a7560759-331f-4243-83c1-b7e49ded50f9.md"
$wsDe.Range("J2").Value = "a7560759-331f-4243-83c1-b7e49ded50f9.b5f762ddd46277abdae7fc8847cc2bcde2ab6a0f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 07:48:23"

$wsDe.Range("I3").Value = "e1f29d49-765f-4abe-8a5d-7a268dab63cd.md"
$wsDe.Range("J3").Value = "e1f29d49-765f-4abe-8a5d-7a268dab63cd.f93c273936f3d89f087d4d48a10478ef87b4e57f.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 07:48:23"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a7560759-331f-4243-83c1-b7e49ded50f9.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlE, "", "", "e1f29d49-765f-4abe-8a5d-7a268dab63cd.md")

$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("I3").Style = "HyperLink"
$wsDe.Range("I2:I3").Font.Underline = $true
$wsDe.Range("I2:I3").Font.Color = 15570276

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# Overview sheet - zh-cn / de-de status columns widen along with the sheets
# above (same shared status text), matching column widths explicitly too.
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Columns.Item(5).ColumnWidth = 29.15
$wsOv.Columns.Item(6).ColumnWidth = 29.15
